$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 172672
$ws.Range("C4").Value = 163463
$ws.Range("C7").Value = 5.33
$ws.Range("C8").Value = 66.05
